# Estadisticos Matutinos 15 Oct
# Update rows 4-7 (groups 3ALCV/5ALCV/5ALCM) on the "1er Parcial" and
# "3er Parcial" sheets with the final grade statistics (Aprobados,
# Reprobados, Por_Apro, Por_Repro, Promedio, Blancos, Por_Blan).

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

# Columns: D=Totales E=Aprobados F=Reprobados G=Por_Apro H=Por_Repro
#          I=Promedio J=Blancos K=Por_Blan
$updates = @{
    4 = @{ E = 22; F = 9;  G = 70.97;  H = 29.03; I = 7.1; J = 0; K = 0 }
    5 = @{ E = 26; F = 8;  G = 76.47;  H = 23.53; I = 7.3; J = 0; K = 0 }
    6 = @{ E = 31; F = 0;  G = 100;    H = 0;     I = 7.7; J = 0; K = 0 }
    7 = @{ E = 34; F = 0;  G = 100;    H = 0;     I = 8.6; J = 0; K = 0 }
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $vals = $updates[$row]
        $ws.Range("E$row").Value = $vals.E
        $ws.Range("F$row").Value = $vals.F
        $ws.Range("G$row").Value = $vals.G
        $ws.Range("H$row").Value = $vals.H
        $ws.Range("I$row").Value = $vals.I
        $ws.Range("J$row").Value = $vals.J
        $ws.Range("K$row").Value = $vals.K
    }
}
